# Remove column A (the per-row index column holding 1, 17) and shift the
# remaining columns (old B:F -> new A:E) one position to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Delete()
